$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = "2025/12/05 15:00"
$ws.Range("B69").Value = "-"
$ws.Range("C69").Value = "-"
$ws.Range("D69").Value = "-"
$ws.Range("E69").Value = "-"
$ws.Range("F69").Value = "-"
$ws.Range("G69").Value = "-"
